$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# The "存款" (deposit) sheet originally had a malformed header row (row 1)
# that just repeated the data values instead of real field names, and it
# was missing the trailing metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) that the other
# sheets (land/building) already carry. Bring it in line with that pattern.

# Extend the header/data formatting (bold + border for row 1, plain for
# row 2) across the new columns before filling them in.
$ws3.Cells.Item(1, 5).Copy()
$ws3.Range("G1:M1").PasteSpecial(-4122)
$ws3.Cells.Item(2, 6).Copy()
$ws3.Range("G2:M2").PasteSpecial(-4122)

# Fix header row 1: B1:D1 become real field names instead of copies of the
# data values, and extend through M1 with the remaining field names.
$ws3.Cells.Item(1, 2).Value = "bank"
$ws3.Cells.Item(1, 3).Value = "deposit_type"
$ws3.Cells.Item(1, 4).Value = "currency"
$ws3.Cells.Item(1, 5).Value = "owner"
$ws3.Cells.Item(1, 6).Value = "total"
$ws3.Cells.Item(1, 7).Value = "property_category"
$ws3.Cells.Item(1, 8).Value = "category"
$ws3.Cells.Item(1, 9).Value = "date"
$ws3.Cells.Item(1, 10).Value = "legislator_name"
$ws3.Cells.Item(1, 11).Value = "legislator_id"
$ws3.Cells.Item(1, 12).Value = "source_file"
$ws3.Cells.Item(1, 13).Value = "index"

# Extend data row 2 with the same metadata already present on the other
# sheets for this legislator/filing.
$ws3.Cells.Item(2, 7).Value = "deposit"
$ws3.Cells.Item(2, 8).Value = "normal"
$ws3.Cells.Item(2, 9).NumberFormat = "@"
$ws3.Cells.Item(2, 9).Value = "2012-04-30"
$ws3.Cells.Item(2, 10).Value = "高金素梅"
$ws3.Cells.Item(2, 11).Value = 926
$ws3.Cells.Item(2, 12).Value = "tmp92521"
$ws3.Cells.Item(2, 13).Value = 45
